# The lab08 report had two pictures whose displayed heights were swapped
# by mistake ("Результат работы функции, шифрующей данные" and "Функция,
# дешифрующая данные"). Fix them by restoring the correct (swapped back)
# heights while keeping every picture's on-page width unchanged.
#
# Word's InlineShape.Height/.Width setters re-lock to the picture's native
# (pixel) aspect ratio as soon as either one is written, which would also
# perturb the width. Using ScaleHeight/ScaleWidth (percentages of the
# picture's native, un-cropped size at 914400/96 = 9525 EMU per pixel)
# lets width and height be set independently, so only cy changes.

$d = $word.ActiveDocument

$EMU_PER_PIXEL = 9525   # 96 dpi: 914400 EMU/inch / 96 px/inch

function Set-InlineShapeExtentEmu($shape, $pixelWidth, $pixelHeight, $targetCxEmu, $targetCyEmu) {
    $nativeCxEmu = $pixelWidth * $EMU_PER_PIXEL
    $nativeCyEmu = $pixelHeight * $EMU_PER_PIXEL

    $shape.ScaleWidth = ($targetCxEmu / $nativeCxEmu) * 100
    $shape.ScaleHeight = ($targetCyEmu / $nativeCyEmu) * 100
}

# InlineShapes(2): "Результат работы функции, шифрующей данные" (image/2.png,
# 1005x650 px) — cy goes back from 3449850 to 2034224 (cx stays 5334000).
$shape2 = $d.InlineShapes.Item(2)
Set-InlineShapeExtentEmu $shape2 1005 650 5334000 2034224

# InlineShapes(3): "Функция, дешифрующая данные" (image/3.png, 1256x479 px)
# — cy goes back from 2034224 to 3449850 (cx stays 5334000).
$shape3 = $d.InlineShapes.Item(3)
Set-InlineShapeExtentEmu $shape3 1256 479 5334000 3449850
